$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at G (shifts email/date_in/date_out from G/H/I to H/I/J)
$ws.Range("G1").EntireColumn.Insert()

# New "phone" column header
$ws.Range("G1").Value = "phone"

# New "phone" column value - force text so the leading zero survives,
# then drop the style back to Normal so the cell keeps the workbook's
# default formatting (matching the other plain text cells).
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "0879128379"
$ws.Range("G2").Style = "Normal"

# New trailing columns after date_out (now J)
$ws.Range("K1").Value = "department_code"
$ws.Range("L1").Value = "division_code"
$ws.Range("M1").Value = "sub_division_code"
$ws.Range("N1").Value = "level_code"
$ws.Range("O1").Value = "position_code"
